$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2022) mirrors column R's formatting row by row.
# Values are first written, then the cell's formatting is copied over from
# the matching cell in column R (same row) via copy / paste-special-formats,
# so each S cell ends up with the same number format / font / border as its
# column-R sibling.

$values = @{
    3  = 2022
    4  = 5.5
    5  = 8.5
    6  = 2.6
    7  = 16.3
    8  = 25.2
    9  = 7.1
    10 = 1.6
    11 = 3.2
    12 = "-"
    13 = 7.5
    14 = 10.5
    15 = 4.5
    16 = 11.4
    17 = 16.1
    18 = 6.6
    19 = 1.2
    20 = 2.1
    21 = 0.3
    22 = 1.5
    23 = 2.9
    24 = 0
    25 = 0.9
    26 = 1.7
    27 = 0.2
    28 = 14.3
    29 = 22.7
    30 = 7.3
    31 = 1.1
    32 = 2.2
    33 = "-"
}

foreach ($row in 3..33) {
    $ws.Range("S$row").Value = $values[$row]
    $ws.Range("R$row").Copy()
    $ws.Range("S$row").PasteSpecial(-4122)
}

# Move / refresh the active selection like the author's Excel session ended up.
$ws.Range("T3").Select() | Out-Null
